$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds product data that was swapped out for a different product in
# this crawl snapshot. Columns A, H and K contain numeric-looking text
# (an id and price-like strings) that must stay text, so their number
# format is forced to Text before the value is written (mirrors typing a
# leading apostrophe in Excel) to avoid Excel auto-coercing them to numbers.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "4947421"

$ws.Range("B3").Value = "Oecoplan Taschentuch Calendula Box"
$ws.Range("C3").Value = "/de/inspiration-geschenke/saisonale-promotionen/gesundheit/oecoplan-taschentuch-calendula-box/p/4947421"
$ws.Range("D3").Value = "80ST"
$ws.Range("E3").Value = 17
$ws.Range("F3").Value = 4

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2.30"

$ws.Range("I3").Value = "0.03/1ST"

$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "0.03"

$ws.Range("M3").Value = "['inspiration-geschenke', 'saisonale-promotionen', 'gesundheit']"
$ws.Range("N3").Value = "Oecoplan Taschentuch Calendula Box 2.30 Schweizer Franken"

# The whole crawl was re-run later the same day, so every row's timestamp
# (column O, rows 2-31) moves from 07:12:23 to 21:00:48.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-30 21:00:48"
}
